$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New firebase_id (string, was numeric) for "mike" (row 5)
$ws.Range("H5").Value = "XJWoEcF4qsToA0NHnKnaIlqBnfO2"

# New alias email for "mike" (row 5), and reset the cell's shading/border to the
# plain default style (it no longer matches the other rows' email style).
$ws.Range("C5").Value = "mike.capstonetest@gmail.com"
$ws.Range("C5").Style = "Normal"

# Widen the firebase_id column (H) to fit the new data
$ws.Columns.Item(8).ColumnWidth = 28.25

# Update the saved selection
$ws.Range("E5").Select() | Out-Null
